{"js": "// Applies the dwq_case_style.docx template fix:\n//  1. Re-join the \"{{case.respondent.name.full()}}{% endif %}...\" text\n//     that had stray proofErr-split runs (no visible text change).\n//  2. Re-join the \"{{case.child.comma_and_list()}}...\" text that had\n//     stray proofErr-split runs (no visible text change).\n//  3. \"IN THE {{court_type}}\" -> \"IN THE {{case.court_type}}\"\n//  4. \"{{court_info.court}}\" -> \"{{case.court_name}}\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Map each paragraph to the replacement that should be applied to it,\n// keyed by a short, unique substring so the logic does not depend on\n// paragraph ordering/index.\nconst replacements = [\n  {\n    match: \"case.respondent.name.full()\",\n    newText:\n      \"{{case.respondent.name.full()}}{% endif %}{% if case.child.number_gathered() > 0 %}\",\n  },\n  {\n    match: \"case.child.comma\",\n    newText:\n      \"{{case.child.comma_and_list()}}{% if case.child.number_gathered()==1 %}, A CHILD{% else %}, CHILDREN{%endif%}{%endif%}\",\n  },\n  {\n    match: \"court_type\",\n    newText: \"IN THE {{case.court_type}}\",\n  },\n  {\n    match: \"court_info.court\",\n    newText: \"{{case.court_name}}\",\n  },\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  for (const repl of replacements) {\n    if (text.indexOf(repl.match) !== -1) {\n      // Replace the whole paragraph's text in one shot so that any\n      // stray proofErr (spelling/grammar) markers splitting the runs\n      // are dropped and the text ends up in a single, clean run.\n      const whole = para.getRange(\"Whole\");\n      whole.insertText(repl.newText, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Applies the dwq_case_style.docx template fix:\n#  1. Re-join the \"{{case.respondent.name.full()}}{% endif %}...\" text\n#     that had stray proofErr-split runs (no visible text change).\n#  2. Re-join the \"{{case.child.comma_and_list()}}...\" text that had\n#     stray proofErr-split runs (no visible text change).\n#  3. \"IN THE {{court_type}}\" -> \"IN THE {{case.court_type}}\"\n#  4. \"{{court_info.court}}\" -> \"{{case.court_name}}\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($doc, $findText, $replaceText) {\n    $rng = $doc.Content\n    # wdFindContinue = 1, wdReplaceAll = 2\n    # MatchCase is left $false because this text is shown in ALL CAPS\n    # (w:caps formatting) but stored in the document in mixed/lower case.\n    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    Write-Host (\"Replace '\" + $findText + \"' -> found=\" + $found)\n}\n\nReplace-DocText $d \"{{case.respondent.name.full()}}{% endif %}{% if case.child.number_gathered() > 0 %}\" \"{{case.respondent.name.full()}}{% endif %}{% if case.child.number_gathered() > 0 %}\"\n\nReplace-DocText $d \"{{case.child.comma_and_list()}}{% if case.child.number_gathered()==1 %}, A CHILD{% else %}, CHILDREN{%endif%}{%endif%}\" \"{{case.child.comma_and_list()}}{% if case.child.number_gathered()==1 %}, A CHILD{% else %}, CHILDREN{%endif%}{%endif%}\"\n\nReplace-DocText $d \"IN THE {{court_type}}\" \"IN THE {{case.court_type}}\"\n\nReplace-DocText $d \"{{court_info.court}}\" \"{{case.court_name}}\"\n"}
